$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (Förändrad) for rows 2 through 7 with the new serial date value
for ($row = 2; $row -le 7; $row++) {
    $ws.Cells.Item($row, 3).Value = 45243
}
